$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing last row (297): high/close values were revised
$ws.Cells.Item(297,4).Value2 = 222.64
$ws.Cells.Item(297,6).Value2 = 222.64

# Extend the data down to row 300, re-using the formatting (date number
# format, border, alignment) already applied to row 297.
$ws.Range("A297:G297").Copy()
$ws.Range("A298:G300").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 298 - 2023-05-01
$ws.Cells.Item(298,1).Value2 = 45047.33333333334
$ws.Cells.Item(298,2).Value2 = "FX_IDC:USDARS"
$ws.Cells.Item(298,3).Value2 = 222.64
$ws.Cells.Item(298,4).Value2 = 239.47
$ws.Cells.Item(298,5).Value2 = 222.64
$ws.Cells.Item(298,6).Value2 = 239.47
$ws.Cells.Item(298,7).Value2 = 0

# Row 299 - 2023-06-01
$ws.Cells.Item(299,1).Value2 = 45078.33333333334
$ws.Cells.Item(299,2).Value2 = "FX_IDC:USDARS"
$ws.Cells.Item(299,3).Value2 = 239.47
$ws.Cells.Item(299,4).Value2 = 256.72
$ws.Cells.Item(299,5).Value2 = 239.215
$ws.Cells.Item(299,6).Value2 = 256.7
$ws.Cells.Item(299,7).Value2 = 0

# Row 300 - 2023-07-03
$ws.Cells.Item(300,1).Value2 = 45110.33333333334
$ws.Cells.Item(300,2).Value2 = "FX_IDC:USDARS"
$ws.Cells.Item(300,3).Value2 = 256.7
$ws.Cells.Item(300,4).Value2 = 261.02
$ws.Cells.Item(300,5).Value2 = 256.7
$ws.Cells.Item(300,6).Value2 = 260.95
$ws.Cells.Item(300,7).Value2 = 0
